$d = $word.ActiveDocument

$replacements = @(
    @("632÷4=158, 0", "356÷2=178, 0"),
    @("184÷7=26, 2", "783÷2=391, 1"),
    @("817÷7=116, 5", "244÷2=122, 0"),
    @("958÷6=159, 4", "178÷6=29, 4"),
    @("912÷6=152, 0", "170÷2=85, 0"),
    @("481÷8=60, 1", "190÷2=95, 0"),
    @("838÷2=419, 0", "945÷2=472, 1"),
    @("262÷5=52, 2", "195÷4=48, 3"),
    @("715÷6=119, 1", "104÷3=34, 2"),
    @("595÷5=119, 0", "334÷9=37, 1"),
    @("457÷6=76, 1", "783÷4=195, 3"),
    @("966÷9=107, 3", "579÷5=115, 4"),
    @("728÷3=242, 2", "718÷9=79, 7"),
    @("555÷3=185, 0", "221÷3=73, 2"),
    @("637÷3=212, 1", "335÷9=37, 2"),
    @("993÷3=331, 0", "139÷9=15, 4"),
    @("974÷9=108, 2", "520÷7=74, 2"),
    @("389÷4=97, 1", "874÷4=218, 2"),
    @("494÷6=82, 2", "699÷3=233, 0"),
    @("668÷8=83, 4", "918÷8=114, 6"),
    @("451÷9=50, 1", "689÷5=137, 4"),
    @("871÷9=96, 7", "243÷5=48, 3"),
    @("708÷8=88, 4", "628÷5=125, 3"),
    @("268÷5=53, 3", "757÷7=108, 1"),
    @("212÷4=53, 0", "337÷7=48, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
